$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "工作表4->工作表4改"
$ws.Range("D13").Select() | Out-Null
